# Updates Sheet1 price (D) and volume-change (E) columns for rows 2-51
# to reflect the refreshed crypto market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "28.203.42"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.792.89"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4508"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +15.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3741"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07546"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.283"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.35%  "
$ws.Range("D16").Value = "1.792.70"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06754"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.351"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").Value = "28.221.22"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.343"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").Value = "1.997.05"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.229"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09407"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.41%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2351"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06311"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02325"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.181"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.480"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.203"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6083"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.793"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.021"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07119"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.158"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
